$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
$ws.Range("A9:B10").ClearContents()

$ws.Range("B2").Select()
